$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values for cryptos list refresh.
$ws.Range("D2").Value = "37.347.97"
$ws.Range("E2").Value = "  +4.07%  "
$ws.Range("D3").Value = "2.039.90"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.98"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.649"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.15"
$ws.Range("E7").Value = "  +9.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.397"
$ws.Range("E9").Value = "  +8.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "59.28"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  +6.58%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.904"
$ws.Range("E13").Value = "  -4.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.21"
$ws.Range("E14").Value = "  +20.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.71"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "2.340.89"
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.67"
$ws.Range("E17").Value = "  +6.26%  "
$ws.Range("D18").Value = "2.047.45"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("D19").Value = "37.262.74"
$ws.Range("E19").Value = "  +3.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.91"
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("D21").Value = "0.0₃0878"
$ws.Range("E21").Value = "  +3.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.45"
$ws.Range("E22").Value = "  +4.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.44"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.35"
$ws.Range("E26").Value = "  +2.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.92"
$ws.Range("E27").Value = "  +2.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.71"
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.92"
$ws.Range("E29").Value = "  +2.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.130"
$ws.Range("E30").Value = "  +33.31%  "
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("E33").Value = "  +4.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0624"
$ws.Range("E34").Value = "  +3.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.62"
$ws.Range("E35").Value = "  +4.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.35"
$ws.Range("E36").Value = "  +10.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  -5.34%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("E40").Value = "  +26.19%  "
$ws.Range("E41").Value = "  +3.77%  "
$ws.Range("E42").Value = "  +9.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.04"
$ws.Range("E43").Value = "  +5.79%  "
$ws.Range("E44").Value = "  +5.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.27"
$ws.Range("E45").Value = "  +4.36%  "
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "94.91"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.80"
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").Value = "1.392.18"
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.90"
$ws.Range("E51").Value = "  +0.20%  "
